$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 3
$ws.Range("B7").Value = 1
$ws.Range("B22").Value = 1
$ws.Range("B47").Value = 3
$ws.Range("B84").Value = 4
$ws.Range("B87").Value = 4
